$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.363865666666667
$ws.Range("H2").Value = 4.091597
$ws.Range("I2").Value = 0.35258381842799
$ws.Range("J2").Value = 0.35258381842799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 42.909214
$ws.Range("N2").Value = 128.727642
$ws.Range("O2").Value = 0.2422627718984814
$ws.Range("P2").Value = 0.2422627718984814
$ws.Range("Q2").Value = 58.52240375825266
$ws.Range("R2").Value = 526.701633824274
$ws.Range("S2").Value = 0.08541793317891572
$ws.Range("T2").Value = 0.08541793317891572
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.363865666666667
$ws.Range("H3").Value = 4.091597
$ws.Range("I3").Value = 0.35258381842799
$ws.Range("J3").Value = 0.35258381842799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 53.36146666666667
$ws.Range("N3").Value = 160.0844
$ws.Range("O3").Value = 0.3012755448569878
$ws.Range("P3").Value = 0.3012755448569878
$ws.Range("Q3").Value = 72.77787230964445
$ws.Range("R3").Value = 655.0008507868001
$ws.Range("S3").Value = 0.1062248820046499
$ws.Range("T3").Value = 0.1062248820046499
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.363865666666667
$ws.Range("H4").Value = 4.091597
$ws.Range("I4").Value = 0.35258381842799
$ws.Range("J4").Value = 0.35258381842799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 39.98186433333333
$ws.Range("N4").Value = 119.945593
$ws.Range("O4").Value = 0.2257351364921847
$ws.Range("P4").Value = 0.2257351364921847
$ws.Range("Q4").Value = 54.52989205355789
$ws.Range("R4").Value = 490.769028482021
$ws.Range("S4").Value = 0.07959055637777798
$ws.Range("T4").Value = 0.07959055637777798
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.363865666666667
$ws.Range("H5").Value = 4.091597
$ws.Range("I5").Value = 0.35258381842799
$ws.Range("J5").Value = 0.35258381842799
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.39145933333334
$ws.Range("N5").Value = 103.174378
$ws.Range("O5").Value = 0.1941720551610951
$ws.Range("P5").Value = 0.1941720551610951
$ws.Range("Q5").Value = 46.90533061129623
$ws.Range("R5").Value = 422.1479755016661
$ws.Range("S5").Value = 0.0684619246407092
$ws.Range("T5").Value = 0.0684619246407092
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.363865666666667
$ws.Range("H6").Value = 4.091597
$ws.Range("I6").Value = 0.35258381842799
$ws.Range("J6").Value = 0.35258381842799
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.474476
$ws.Range("N6").Value = 19.423428
$ws.Range("O6").Value = 0.03655449159125106
$ws.Range("P6").Value = 0.03655449159125106
$ws.Range("Q6").Value = 8.830315526057333
$ws.Range("R6").Value = 79.47283973451601
$ws.Range("S6").Value = 0.01288852222593715
$ws.Range("T6").Value = 0.01288852222593715
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.247734666666667
$ws.Range("H7").Value = 3.743204
$ws.Range("I7").Value = 0.3225618650798028
$ws.Range("J7").Value = 0.3225618650798028
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.909214
$ws.Range("N7").Value = 128.727642
$ws.Range("O7").Value = 0.2422627718984814
$ws.Range("P7").Value = 0.2422627718984814
$ws.Range("Q7").Value = 53.53931382721868
$ws.Range("R7").Value = 481.853824444968
$ws.Range("S7").Value = 0.07814473154297701
$ws.Range("T7").Value = 0.07814473154297699
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.247734666666667
$ws.Range("H8").Value = 3.743204
$ws.Range("I8").Value = 0.3225618650798028
$ws.Range("J8").Value = 0.3225618650798028
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 53.36146666666667
$ws.Range("N8").Value = 160.0844
$ws.Range("O8").Value = 0.3012755448569878
$ws.Range("P8").Value = 0.3012755448569878
$ws.Range("Q8").Value = 66.5809518241778
$ws.Range("R8").Value = 599.2285664176002
$ws.Range("S8").Value = 0.09718000165200377
$ws.Range("T8").Value = 0.09718000165200376
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.247734666666667
$ws.Range("H9").Value = 3.743204
$ws.Range("I9").Value = 0.3225618650798028
$ws.Range("J9").Value = 0.3225618650798028
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 39.98186433333333
$ws.Range("N9").Value = 119.945593
$ws.Range("O9").Value = 0.2257351364921847
$ws.Range("P9").Value = 0.2257351364921847
$ws.Range("Q9").Value = 49.88675816666356
$ws.Range("R9").Value = 448.9808234999721
$ws.Range("S9").Value = 0.07281354664096296
$ws.Range("T9").Value = 0.07281354664096294
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.247734666666667
$ws.Range("H10").Value = 3.743204
$ws.Range("I10").Value = 0.3225618650798028
$ws.Range("J10").Value = 0.3225618650798028
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 34.39145933333334
$ws.Range("N10").Value = 103.174378
$ws.Range("O10").Value = 0.1941720551610951
$ws.Range("P10").Value = 0.1941720551610951
$ws.Range("Q10").Value = 42.9114160474569
$ws.Range("R10").Value = 386.2027444271121
$ws.Range("S10").Value = 0.06263250025914119
$ws.Range("T10").Value = 0.06263250025914117
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.247734666666667
$ws.Range("H11").Value = 3.743204
$ws.Range("I11").Value = 0.3225618650798028
$ws.Range("J11").Value = 0.3225618650798028
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 6.474476
$ws.Range("N11").Value = 19.423428
$ws.Range("O11").Value = 0.03655449159125106
$ws.Range("P11").Value = 0.03655449159125106
$ws.Range("Q11").Value = 8.078428153701335
$ws.Range("R11").Value = 72.70585338331202
$ws.Range("S11").Value = 0.01179108498471791
$ws.Range("T11").Value = 0.01179108498471791
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.808894
$ws.Range("H12").Value = 2.426682
$ws.Range("I12").Value = 0.2091136555409713
$ws.Range("J12").Value = 0.2091136555409713
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 42.909214
$ws.Range("N12").Value = 128.727642
$ws.Range("O12").Value = 0.2422627718984814
$ws.Range("P12").Value = 0.2422627718984814
$ws.Range("Q12").Value = 34.709005749316
$ws.Range("R12").Value = 312.381051743844
$ws.Range("S12").Value = 0.05066045383317994
$ws.Range("T12").Value = 0.05066045383317994
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.808894
$ws.Range("H13").Value = 2.426682
$ws.Range("I13").Value = 0.2091136555409713
$ws.Range("J13").Value = 0.2091136555409713
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 53.36146666666667
$ws.Range("N13").Value = 160.0844
$ws.Range("O13").Value = 0.3012755448569878
$ws.Range("P13").Value = 0.3012755448569878
$ws.Range("Q13").Value = 43.16377021786667
$ws.Range("R13").Value = 388.4739319608
$ws.Range("S13").Value = 0.06300083051014259
$ws.Range("T13").Value = 0.06300083051014259
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.808894
$ws.Range("H14").Value = 2.426682
$ws.Range("I14").Value = 0.2091136555409713
$ws.Range("J14").Value = 0.2091136555409713
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 39.98186433333333
$ws.Range("N14").Value = 119.945593
$ws.Range("O14").Value = 0.2257351364921847
$ws.Range("P14").Value = 0.2257351364921847
$ws.Range("Q14").Value = 32.34109016804734
$ws.Range("R14").Value = 291.069811512426
$ws.Range("S14").Value = 0.04720429957592086
$ws.Range("T14").Value = 0.04720429957592085
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.808894
$ws.Range("H15").Value = 2.426682
$ws.Range("I15").Value = 0.2091136555409713
$ws.Range("J15").Value = 0.2091136555409713
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 34.39145933333334
$ws.Range("N15").Value = 103.174378
$ws.Range("O15").Value = 0.1941720551610951
$ws.Range("P15").Value = 0.1941720551610951
$ws.Range("Q15").Value = 27.81904510597734
$ws.Range("R15").Value = 250.371405953796
$ws.Range("S15").Value = 0.04060402825863972
$ws.Range("T15").Value = 0.04060402825863971
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.808894
$ws.Range("H16").Value = 2.426682
$ws.Range("I16").Value = 0.2091136555409713
$ws.Range("J16").Value = 0.2091136555409713
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.474476
$ws.Range("N16").Value = 19.423428
$ws.Range("O16").Value = 0.03655449159125106
$ws.Range("P16").Value = 0.03655449159125106
$ws.Range("Q16").Value = 5.237164789544
$ws.Range("R16").Value = 47.134483105896
$ws.Range("S16").Value = 0.007644043363088207
$ws.Range("T16").Value = 0.007644043363088207
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4477083333333334
$ws.Range("H17").Value = 1.343125
$ws.Range("I17").Value = 0.1157406609512359
$ws.Range("J17").Value = 0.1157406609512359
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 42.909214
$ws.Range("N17").Value = 128.727642
$ws.Range("O17").Value = 0.2422627718984814
$ws.Range("P17").Value = 0.2422627718984814
$ws.Range("Q17").Value = 19.21081268458333
$ws.Range("R17").Value = 172.89731416125
$ws.Range("S17").Value = 0.02803965334340875
$ws.Range("T17").Value = 0.02803965334340875
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.4477083333333334
$ws.Range("H18").Value = 1.343125
$ws.Range("I18").Value = 0.1157406609512359
$ws.Range("J18").Value = 0.1157406609512359
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 53.36146666666667
$ws.Range("N18").Value = 160.0844
$ws.Range("O18").Value = 0.3012755448569878
$ws.Range("P18").Value = 0.3012755448569878
$ws.Range("Q18").Value = 23.89037330555556
$ws.Range("R18").Value = 215.01335975
$ws.Range("S18").Value = 0.03486983069019149
$ws.Range("T18").Value = 0.03486983069019149
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.4477083333333334
$ws.Range("H19").Value = 1.343125
$ws.Range("I19").Value = 0.1157406609512359
$ws.Range("J19").Value = 0.1157406609512359
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 39.98186433333333
$ws.Range("N19").Value = 119.945593
$ws.Range("O19").Value = 0.2257351364921847
$ws.Range("P19").Value = 0.2257351364921847
$ws.Range("Q19").Value = 17.90021384423611
$ws.Range("R19").Value = 161.101924598125
$ws.Range("S19").Value = 0.02612673389752292
$ws.Range("T19").Value = 0.02612673389752291
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 0.4477083333333334
$ws.Range("H20").Value = 1.343125
$ws.Range("I20").Value = 0.1157406609512359
$ws.Range("J20").Value = 0.1157406609512359
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 34.39145933333334
$ws.Range("N20").Value = 103.174378
$ws.Range("O20").Value = 0.1941720551610951
$ws.Range("P20").Value = 0.1941720551610951
$ws.Range("Q20").Value = 15.39734293902778
$ws.Range("R20").Value = 138.57608645125
$ws.Range("S20").Value = 0.02247360200260499
$ws.Range("T20").Value = 0.02247360200260498
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 0.4477083333333334
$ws.Range("H21").Value = 1.343125
$ws.Range("I21").Value = 0.1157406609512359
$ws.Range("J21").Value = 0.1157406609512359
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 6.474476
$ws.Range("N21").Value = 19.423428
$ws.Range("O21").Value = 0.03655449159125106
$ws.Range("P21").Value = 0.03655449159125106
$ws.Range("Q21").Value = 2.898676859166667
$ws.Range("R21").Value = 26.0880917325
$ws.Range("S21").Value = 0.004230841017507795
$ws.Range("T21").Value = 0.004230841017507795
